# Update annotations for parisk
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: politeness_score becomes a real number (was text "3"),
# and polite_expressions is cleared (was text "nan") to an empty string,
# matching the blank "polite_expressions" cells used elsewhere in the sheet.
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "'"
$ws.Range("C11").Style = "Normal"

# New row 12: another parisk annotation.
$ws.Range("A12").Value = "parisk"
$ws.Range("B12").Value = "'2"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "nan"
$ws.Range("D12").Value = "DFT"
$ws.Range("E12").Value = "EXP"
$ws.Range("F12").Value = "afe80f3f-3501-40b4-a3d0-1ad1f86c76ec"
$ws.Range("G12").Value = "r1BRfhiab_annotated.xlsx"
$ws.Range("H12").Value = "Not too surprisingly, the standard multiclass losses do not have the desired property, however approaches that reduce multi-class to binary classification at training time do, namely unnormalized models with penalized log Z (self-normalization), the NCE approach, as well as (the natural in the proposed setting) binary classification loss."
